$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data: 3 omitted census tracts (plus the 102.02 correction) for Albemarle ---
$rows = @(
    @(" Albemarle", "Census Tract 102.02", 85.876349582455006),
    @(" Albemarle", "Census Tract 109.01", 81.601332669245451),
    @(" Albemarle", "Census Tract 109.02", 81.272027683090585),
    @(" Albemarle", "Census Tract 109.03", 85.150992842316668)
)

$startRow = 48
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

# --- Match formatting of the preceding data block (rows 36/37 alternate fill+number style) ---
$ws.Range("A36:B37").Copy() | Out-Null
$ws.Range("A48:B49").PasteSpecial(-4122) | Out-Null
$ws.Range("A50:B51").PasteSpecial(-4122) | Out-Null

$ws.Range("C36:C37").Copy() | Out-Null
$ws.Range("C48:C49").PasteSpecial(-4122) | Out-Null
$ws.Range("C50:C51").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Move the view/selection to the newly-added data, like the author scrolling down ---
$ws.Range("B51").Select() | Out-Null
